# Commit: "update Pablo 18 11"
#
# 1) Generators_AssetData_Existing: remove the "no fill" highlight that was
#    applied down column F (F1:F17), and clear the stray formatted-but-empty
#    cell F18 so the trailing blank row disappears entirely.
# 2) Fuel_Cost_Absolute: remove the same "no fill" highlight from B4.
# 3) Python_Gen_E_Data: add a new "Capacity" column (D) that mirrors the
#    P_max column (F) from Generators_AssetData_Existing, widen column B a
#    touch, and make this the active/selected sheet.
# 4) Update the selection left on Generators_AssetData_Existing now that it
#    is no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- Sheet: Generators_AssetData_Existing -------------------------------
$wsGen = $wb.Worksheets.Item("Generators_AssetData_Existing")

# Drop the fill/highlight that used to sit on the P_max (F) column.
$wsGen.Range("F1:F17").Interior.Pattern = -4142

# F18 was an empty, but still-formatted, leftover cell -- clear it out
# completely so the row disappears from the sheet.
$wsGen.Range("F18").Clear()

# --- Sheet: Fuel_Cost_Absolute -------------------------------------------
$wsFuel = $wb.Worksheets.Item("Fuel_Cost_Absolute")
$wsFuel.Range("B4").Interior.Pattern = -4142

# --- Sheet: Python_Gen_E_Data --------------------------------------------
$wsPy = $wb.Worksheets.Item("Python_Gen_E_Data")

$wsPy.Range("D1").Value = "Capacity"
for ($r = 2; $r -le 17; $r++) {
    $wsPy.Range("D$r").Formula = "=Generators_AssetData_Existing!F$r"
}

$wsPy.Columns.Item(2).ColumnWidth = 11.33

# --- Selection / active-sheet bookkeeping --------------------------------
$wsGen.Range("H19").Select()

$wsPy.Activate()
$wsPy.Range("D2:D17").Select()
